# Fix report 2 template: switch the dissertation table to a fixed
# layout and resize its columns.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Equivalent to adding <w:tblLayout w:type="fixed"/> to the table
# properties (Word auto-layout -> fixed layout).
$t.AllowAutoFit = $false

# New column widths in twips (dxa), converted to points (Word's COM
# Width property uses points, 1 pt = 20 dxa).
$newWidthsDxa = @(568, 1701, 1418, 1108, 1585, 1276, 1275, 1240)

for ($i = 1; $i -le $t.Columns.Count; $i++) {
    $t.Columns.Item($i).Width = $newWidthsDxa[$i - 1] / 20.0
}
